# Added try-catch in SQLHandler.
#
# Replaces the seven "Name|Name" placeholder values (stored as shared
# strings used throughout column B) with their corresponding e-mail
# addresses, drops the custom width override on column B, and updates
# the sheet's selection to span the whole of column B with the active
# cell positioned further down (matching the scrolled view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# Replace each shared "name" string with its e-mail counterpart, in the
# same order the strings originally appear in the shared-string table so
# the table positions line up with the intended result.
$used.Replace("Alex|Krylov", "slate@list.ru")
$used.Replace("Petrov|Vladimir", "vlad@list.ru")
$used.Replace("Popov|Sergei", "ser@list.ru")
$used.Replace("Ivanov|Albert", "albiv@list.ru")
$used.Replace("Serov|Valera", "vals@list.ru")
$used.Replace("Ponov|Pavel", "popov@list.ru")
$used.Replace("Okolov|Seva", "seva@list.ru")

# Drop the custom column width that had been applied to column B.
$ws.Columns.Item(2).ClearFormats()

# Select the entire column B (matches the sheet's saved selection of the
# whole column, scrolled down so row 694 is in view).
$ws.Columns.Item(2).Select()
